$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Rename label in column A, row 6 (shared string "DN_ConsSD" -> "DL_ConsSD")
$ws.Range("A6").Value = "DL_ConsSD"

# Row 3 (XG_ConsSD) - small precision updates
$ws.Range("C3").Value = 0.008564035773174005
$ws.Range("D3").Value = 0.004427547437360639
$ws.Range("E3").Value = 0.01634454154301983
$ws.Range("G3").Value = 0.01496273640566521
$ws.Range("I3").Value = 0.03256071053709742

# Row 6 (DL_ConsSD) - updated metric values
$ws.Range("B6").Value = 0.9982786914031656
$ws.Range("C6").Value = 0.03161798606353385
$ws.Range("D6").Value = 0.01468101875194277
$ws.Range("E6").Value = 0.04528331523957421
$ws.Range("F6").Value = 0.9818873422160905
$ws.Range("G6").Value = 0.1128581830842473
$ws.Range("H6").Value = 0.06060163320539391
$ws.Range("I6").Value = 0.146892606281703
